# Fix random walk estimation and init edit.
# On the "meta_data_names" sheet, columns F:G (rows 12-31) hold a small
# glossary of fleet_control field names + descriptions. A new entry,
# "Age_max_selected" (age at which selectivity = 1; negative normalizes by
# max), is inserted right after "Age_first_selected" (row 11), which pushes
# every following F:G glossary entry down by one row. The last entry ("Sex")
# moves from row 30 to row 31, so row 30 loses its F:G pair (but keeps/gains
# its normal A:E "fleet_control" row content).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-FG($row, $fText, $gText, $style) {
    $fCell = $ws.Range("F$row")
    $gCell = $ws.Range("G$row")
    if ($null -eq $fText) {
        $fCell.ClearContents()
        $gCell.ClearContents()
        return
    }
    $fCell.Value2 = $fText
    $gCell.Value2 = $gText
    if ($style -eq 2) {
        $fCell.ClearFormats()
        $fCell.HorizontalAlignment = -4131
    } elseif ($style -eq 3) {
        $fCell.ClearFormats()
        $fCell.HorizontalAlignment = -4131
        $fCell.Font.Color = 2236962
    } else {
        $fCell.ClearFormats()
    }
}

# New entry inserted at row 12 (no special style, like its F11 neighbor)
Set-FG 12 "Age_max_selected" "Age at which selectivity = 1. If NA, it will not normalize selectivity. If < 0, will normalize selectivity by the max." 0
Set-FG 13 "Comp_loglike" "Composition data distribution (0 = multinomial; 1 = dirichlet-multinomial)" 0
Set-FG 14 "weight1_Numbers2" "Is the observation in weight (kg) set as 1, if the observation is in numbers caught, set as 2" 3
Set-FG 15 "Weight_index" "Weight-at-age (wt) index to use for calculation of derived quantities" 0
Set-FG 16 "Age_transition_index" "Age transition matrix (e.g. growth trajectory) index to use for derived quantities to convert age to length" 0
Set-FG 17 "Q_index" "index to use if catchability coefficients are to be set the same" 0
Set-FG 18 "Estimate_q" "Estimate catchability? (0 = fixed at prior; - 1 = Estimate single parameter; - 2 = Estimate single parameter with prior; - 3 = Estimate analytical q  from Ludwig and Walters 1994;  - 4 = Estimate power equation; - 5 - Linear equation log(q_y) = q_mu + beta * index_y)" 0
Set-FG 19 "Q_prior" "Starting value or fixed value for catchability" 2
Set-FG 20 "Q_sd_prior" "Variance of q prior: dnorm (log_q, log_q_prior, q_sd_prior)" 0
Set-FG 21 "Time_varying_q" "Wether a time-varying q should be estimated. 0 = no, 1 = penalized deviate, 3 = time blocks with no penalty; 4 = random walk from mean following Dorn 2018 (dnorm(q_y - q_y-1, 0, sigma). If Estimate_q = 5, this determines the environmental index to be used in the equation log(q_y) = q_mu + beta * index_y. If `"random_q`" is selected in fit_mod, penalized deviates (1) and random walk parameters (4) will be treated as random effects." 0
Set-FG 22 "Time_varying_q_sd_prior" "The sd to use for the random walk of time varying q if set to 1" 0
Set-FG 23 "Estimate_survey_sd" "Estimate survey variance (0 = use CV from index_data, 1 = yes, 2 = analytically estimate following (Ludwig and Walters 1994)" 2
Set-FG 24 "Survey_sd_prior" "Starting value to be used if Estimate_sigma_index = 1" 2
Set-FG 25 "Estimate_catch_sd" "Estimate fishery variance (0 = use CV from index_data, 1 = yes, 2 = analytically estimate following (Ludwig and Walters 1994)" 2
Set-FG 26 "Catch_sd_prior" "Starting value to be used if Estimate_sigma_catch = 1" 2
Set-FG 27 "Comp_weights" "Composition weights to be used for multinomial likelihood. These are multiplied. After running model, these will update to McAllister & Ianelli 1997 weights using the harmonic mean." 2
Set-FG 28 "Catch_units" "Units used for survey: 1 = kg; 2 = numbers" 0
Set-FG 29 "proj_F_prop" "The proportion of future fishing mortality assigned to this fleet" 2
# Row 30 no longer has an F:G glossary pair
Set-FG 30 $null $null $null
Set-FG 31 "Sex" "sex codes: 0=combined; 1=use female only; 2=use male only; 3 = joint female and male" 0

# Row 30's A:D content (propF row) is untouched, but it now also carries the
# "fleet_control" marker in column E, matching the other rows in this block.
$ws.Range("E30").Value2 = "fleet_control"

# Leave the same cell selected/active as in the edited workbook
$ws.Range("D13").Select()
